$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.659.62"
$ws.Range("D3").Value = "'1.953.08"
$ws.Range("D4").Value = "'0.9995"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'247.42"
$ws.Range("E5").Value = "  +1.11%  "
$ws.Range("D6").Value = "'0.9991"
$ws.Range("D7").Value = "'0.4814"
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  +1.73%  "
$ws.Range("D10").Value = "'0.06800"
$ws.Range("E10").Value = "  +1.25%  "
$ws.Range("D11").Value = "'112.05"
$ws.Range("E11").Value = "  +1.47%  "
$ws.Range("D12").Value = "'19.45"
$ws.Range("E12").Value = "  +2.18%  "
$ws.Range("D13").Value = "'1.966.24"
$ws.Range("E13").Value = "  +2.85%  "
$ws.Range("D14").Value = "'0.07690"
$ws.Range("E14").Value = "  +1.91%  "
$ws.Range("E15").Value = "  +4.23%  "
$ws.Range("D16").Value = "'0.6853"
$ws.Range("E16").Value = "  +2.21%  "
$ws.Range("D17").Value = "'293.85"
$ws.Range("E17").Value = "  +1.77%  "
$ws.Range("D18").Value = "'30.667.01"
$ws.Range("E18").Value = "  +0.65%  "
$ws.Range("B19").Value = "BitDAO"
$ws.Range("C19").Value = "https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit"
$ws.Range("D19").Value = "'0.4998"
$ws.Range("E19").Value = "  +18.05%  "
$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").Value = "'13.23"
$ws.Range("E20").Value = "  +3.11%  "
$ws.Range("D21").Value = "'2.225.39"
$ws.Range("E21").Value = "  +2.86%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "'5.650"
$ws.Range("E22").Value = "  +3.08%  "
$ws.Range("B23").Value = "ShibaInu"
$ws.Range("C23").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D23").Value = "'0.000007670"
$ws.Range("E23").Value = "  +1.20%  "
$ws.Range("B24").Value = "Dai"
$ws.Range("C24").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D24").Value = "'1.000"
$ws.Range("E24").Value = "  +0.14%  "
$ws.Range("B25").Value = "BinanceUSD"
$ws.Range("C25").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D25").Value = "'0.9995"
$ws.Range("E25").Value = "  +0.08%  "
$ws.Range("B26").Value = "Chainlink"
$ws.Range("C26").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D26").Value = "'6.606"
$ws.Range("E26").Value = "  +3.14%  "
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").Value = "'9.755"
$ws.Range("E27").Value = "  +3.29%  "
$ws.Range("B28").Value = "Monero"
$ws.Range("C28").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D28").Value = "'168.90"
$ws.Range("E28").Value = "  +2.84%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "'20.24"
$ws.Range("E29").Value = "  -0.48%  "
$ws.Range("B30").Value = "LidoDAOToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D30").Value = "'2.192"
$ws.Range("E30").Value = "  +3.68%  "
$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").Value = "'0.1087"
$ws.Range("E31").Value = "  +3.22%  "
$ws.Range("B32").Value = "Toncoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D32").Value = "'1.435"
$ws.Range("E32").Value = "  +2.28%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "'4.635"
$ws.Range("E33").Value = "  +14.83%  "
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").Value = "'4.425"
$ws.Range("E34").Value = "  +6.27%  "
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "'0.05059"
$ws.Range("E35").Value = "  +1.68%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "'0.7753"
$ws.Range("E36").Value = "  +6.40%  "
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").Value = "'1.164"
$ws.Range("E37").Value = "  +2.82%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.02074"
$ws.Range("E38").Value = "  +2.02%  "
$ws.Range("B39").Value = "HuobiToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D39").Value = "'2.730"
$ws.Range("E39").Value = "  +0.37%  "
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value = "'2.698"
$ws.Range("E40").Value = "  +1.16%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").Value = "'2.059"
$ws.Range("E41").Value = "  +2.19%  "
$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").Value = "'110.83"
$ws.Range("E42").Value = "  +0.29%  "
$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").Value = "'0.4459"
$ws.Range("E43").Value = "  +0.47%  "
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").Value = "'0.8726"
$ws.Range("E44").Value = "  +0.87%  "
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").Value = "'5.988"
$ws.Range("E45").Value = "  +3.53%  "
$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").Value = "'0.9993"
$ws.Range("E46").Value = "  +0.06%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "'69.48"
$ws.Range("E47").Value = "  +1.96%  "
$ws.Range("B48").Value = "Aptos"
$ws.Range("C48").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D48").Value = "'7.385"
$ws.Range("E48").Value = "  +0.85%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'9.337"
$ws.Range("E49").Value = "  +1.00%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "'0.1253"
$ws.Range("E50").Value = "  +1.04%  "
$ws.Range("B51").Value = "BitcoinSV"
$ws.Range("C51").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D51").Value = "'47.94"
$ws.Range("E51").Value = "  -2.36%  "
